# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update "Conversión del día" text with new Binance rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.58 = 41408.05 pesos`n✅ 41408.05 pesos = 9.56 = 958.75 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- tasas!N10/O10/N12/O12: update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 104.4
$ws2.Range("O10").Value = 4323
$ws2.Range("N12").Value = 4333
$ws2.Range("O12").Value = 100.325
